# Update cryptos list values per latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.180.44'
$ws.Range("E2").Value = '  -0.69%  '
$ws.Range("D3").Value = '3.424.52'
$ws.Range("E3").Value = '  +0.12%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''571.90'
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").Value = '''161.42'
$ws.Range("E6").Value = '  +2.23%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.425.24'
$ws.Range("E8").Value = '  +0.17%  '
$ws.Range("D9").Value = '''0.555'
$ws.Range("E9").Value = '  -8.12%  '
$ws.Range("D10").Value = '''7.29'
$ws.Range("E10").Value = '  +1.53%  '
$ws.Range("E11").Value = '  -1.85%  '
$ws.Range("E12").Value = '  -3.12%  '
$ws.Range("D13").Value = '4.017.34'
$ws.Range("E13").Value = '  +0.23%  '
$ws.Range("D15").Value = '''27.09'
$ws.Range("E15").Value = '  -1.86%  '
$ws.Range("E16").Value = '  -6.50%  '
$ws.Range("D17").Value = '64.198.42'
$ws.Range("E17").Value = '  -0.75%  '
$ws.Range("D18").Value = '3.436.80'
$ws.Range("E18").Value = '  +1.06%  '
$ws.Range("E19").Value = '  -3.90%  '
$ws.Range("D20").Value = '''13.59'
$ws.Range("E20").Value = '  -1.34%  '
$ws.Range("D21").Value = '''378.00'
$ws.Range("E21").Value = '  -0.46%  '
$ws.Range("E22").Value = '  -1.79%  '
$ws.Range("D24").Value = '''71.48'
$ws.Range("E24").Value = '  -0.59%  '
$ws.Range("D25").Value = '''0.518'
$ws.Range("E25").Value = '  -5.30%  '
$ws.Range("E26").Value = '  -1.29%  '
$ws.Range("E27").Value = '  -4.02%  '
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  -0.17%  '
$ws.Range("D30").Value = '''6.02'
$ws.Range("E30").Value = '  -2.14%  '
$ws.Range("E31").Value = '  -3.93%  '
$ws.Range("E32").Value = '  +0.55%  '
$ws.Range("D33").Value = '''22.94'
$ws.Range("E33").Value = '  -1.12%  '
$ws.Range("D34").Value = '''7.07'
$ws.Range("E34").Value = '  +0.43%  '
$ws.Range("E35").Value = '  -3.84%  '
$ws.Range("D36").Value = '''159.61'
$ws.Range("E36").Value = '  -0.72%  '
$ws.Range("E37").Value = '  +11.53%  '
$ws.Range("E38").Value = '  -4.27%  '
$ws.Range("D39").Value = '2.810.35'
$ws.Range("E39").Value = '  -2.44%  '
$ws.Range("D40").Value = '''0.0727'
$ws.Range("E40").Value = '  -3.55%  '
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").Value = '''25.80'
$ws.Range("E41").Value = '  -1.90%  '
$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").Value = '''42.99'
$ws.Range("E42").Value = '  -0.15%  '
$ws.Range("D43").Value = '''6.47'
$ws.Range("E43").Value = '  -3.21%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").Value = '''4.44'
$ws.Range("E44").Value = '  -2.58%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '''26.12'
$ws.Range("E45").Value = '  +1.16%  '
$ws.Range("D46").Value = '''0.0305'
$ws.Range("E46").Value = '  -3.41%  '
$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D47").Value = '''2.41'
$ws.Range("E47").Value = '  +8.71%  '
$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").Value = '''336.78'
$ws.Range("E48").Value = '  +6.87%  '
$ws.Range("D49").Value = '''1.05'
$ws.Range("E49").Value = '  -0.97%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = '''6.32'
$ws.Range("E50").Value = '  -3.06%  '
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").Value = '''0.103'
$ws.Range("E51").Value = '  -3.79%  '
